$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change MP time limit status: "OPTIMAL" -> "TIME_LIMIT" (shared across all Status cells E2:E11)
$ws.Range("E2:E11").Value = "TIME_LIMIT"

# Corrected fixed recourse data: objective (B), gap (C), solve time (D) for rows 2-11
$data = @(
    @(-1442.4032829918265, 9.611728021278704, 5710.12178695),
    @(-1440.2106322599902, 9.110605400142823, 5737.917843119),
    @(-1334.489771790431, 19.194195800119036, 5541.354182354),
    @(-1439.5884369399998, 9.681649625434384, 5849.72208816),
    @(-1449.8592891489982, 8.45484907740723, 5695.82920966),
    @(-1429.6423532131853, 9.77376211826315, 5651.239159975),
    @(-1276.6756221850605, 24.400064142381794, 5840.649874044),
    @(-1315.54388657078, 19.11427061444898, 5853.535549383),
    @(-1299.8603260651262, 19.87148461062062, 5906.744287678),
    @(-1265.0477944423028, 25.690724567472962, 5625.023825159)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $row++
}
